# Auto-generated Excel COM-interop script to apply cryptos.xlsx data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.044.91"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "1.902.24"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.04"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4644"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4118"
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.58"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07988"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.83"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "1.902.17"
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.939"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.090"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001033"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06584"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "29.108.01"
$ws.Range("E22").Value = "  +2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.439"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.218"
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("D26").Value = "2.127.93"
$ws.Range("E26").Value = "  +2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.41"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.72"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.127"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.435"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.30"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9829"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09410"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.595"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.303"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06097"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02243"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.350"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.176"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5805"
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1822"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.263"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.309"
$ws.Range("E46").Value = "  +11.94%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.11"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5511"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.914"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07054"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.64"
$ws.Range("E51").Value = "  +18.10%  "
